$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.750268
$ws.Range("H2").Value = 59.250804
$ws.Range("I2").Value = 0.2138888518073023
$ws.Range("J2").Value = 0.2138888518073023
$ws.Range("M2").Value = 5.256051666666667
$ws.Range("N2").Value = 15.768155
$ws.Range("O2").Value = 0.2751978571638378
$ws.Range("P2").Value = 0.2751978571638379
$ws.Range("Q2").Value = 103.8084290385133
$ws.Range("R2").Value = 934.27586134662
$ws.Range("S2").Value = 0.05886175368860326
$ws.Range("T2").Value = 0.05886175368860327
$ws.Range("G3").Value = 19.750268
$ws.Range("H3").Value = 59.250804
$ws.Range("I3").Value = 0.2138888518073023
$ws.Range("J3").Value = 0.2138888518073023
$ws.Range("O3").Value = 0.1142283101567343
$ws.Range("P3").Value = 0.1142283101567343
$ws.Range("Q3").Value = 43.08849476991066
$ws.Range("R3").Value = 387.7964529291959
$ws.Range("S3").Value = 0.02443216210331232
$ws.Range("T3").Value = 0.02443216210331232
$ws.Range("G4").Value = 19.750268
$ws.Range("H4").Value = 59.250804
$ws.Range("I4").Value = 0.2138888518073023
$ws.Range("J4").Value = 0.2138888518073023
$ws.Range("M4").Value = 11.66145566666667
$ws.Range("N4").Value = 34.98436700000001
$ws.Range("O4").Value = 0.6105738326794278
$ws.Range("P4").Value = 0.6105738326794279
$ws.Range("Q4").Value = 230.3168746867854
$ws.Range("R4").Value = 2072.851872181068
$ws.Range("S4").Value = 0.1305949360153867
$ws.Range("T4").Value = 0.1305949360153868
$ws.Range("I5").Value = 0.6395228081370402
$ws.Range("J5").Value = 0.6395228081370402
$ws.Range("M5").Value = 5.256051666666667
$ws.Range("N5").Value = 15.768155
$ws.Range("O5").Value = 0.2751978571638378
$ws.Range("P5").Value = 0.2751978571638379
$ws.Range("Q5").Value = 310.3848446800545
$ws.Range("R5").Value = 2793.46360212049
$ws.Range("S5").Value = 0.1759953064067137
$ws.Range("T5").Value = 0.1759953064067137
$ws.Range("I6").Value = 0.6395228081370402
$ws.Range("J6").Value = 0.6395228081370402
$ws.Range("O6").Value = 0.1142283101567343
$ws.Range("P6").Value = 0.1142283101567343
$ws.Range("S6").Value = 0.07305160968018354
$ws.Range("T6").Value = 0.07305160968018354
$ws.Range("I7").Value = 0.6395228081370402
$ws.Range("J7").Value = 0.6395228081370402
$ws.Range("M7").Value = 11.66145566666667
$ws.Range("N7").Value = 34.98436700000001
$ws.Range("O7").Value = 0.6105738326794278
$ws.Range("P7").Value = 0.6105738326794279
$ws.Range("Q7").Value = 688.6422233625319
$ws.Range("R7").Value = 6197.780010262787
$ws.Range("S7").Value = 0.390475892050143
$ws.Range("T7").Value = 0.390475892050143
$ws.Range("G8").Value = 13.53581066666667
$ws.Range("H8").Value = 40.607432
$ws.Range("I8").Value = 0.1465883400556574
$ws.Range("J8").Value = 0.1465883400556574
$ws.Range("M8").Value = 5.256051666666667
$ws.Range("N8").Value = 15.768155
$ws.Range("O8").Value = 0.2751978571638378
$ws.Range("P8").Value = 0.2751978571638379
$ws.Range("Q8").Value = 71.14492021421779
$ws.Range("R8").Value = 640.3042819279601
$ws.Range("S8").Value = 0.0403407970685209
$ws.Range("T8").Value = 0.04034079706852091
$ws.Range("G9").Value = 13.53581066666667
$ws.Range("H9").Value = 40.607432
$ws.Range("I9").Value = 0.1465883400556574
$ws.Range("J9").Value = 0.1465883400556574
$ws.Range("O9").Value = 0.1142283101567343
$ws.Range("P9").Value = 0.1142283101567343
$ws.Range("Q9").Value = 29.53062242584089
$ws.Range("R9").Value = 265.775601832568
$ws.Range("S9").Value = 0.01674453837323848
$ws.Range("T9").Value = 0.01674453837323848
$ws.Range("G10").Value = 13.53581066666667
$ws.Range("H10").Value = 40.607432
$ws.Range("I10").Value = 0.1465883400556574
$ws.Range("J10").Value = 0.1465883400556574
$ws.Range("M10").Value = 11.66145566666667
$ws.Range("N10").Value = 34.98436700000001
$ws.Range("O10").Value = 0.6105738326794278
$ws.Range("P10").Value = 0.6105738326794279
$ws.Range("Q10").Value = 157.8472560017271
$ws.Range("R10").Value = 1420.625304015544
$ws.Range("S10").Value = 0.08950300461389804
$ws.Range("T10").Value = 0.08950300461389805
